# Add the new "alarm identification" rows (9-12) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append below the existing table (A1:F8).
# Every value in the source workbook is stored as text (even the
# numeric-looking ones), so we force a text number format on each
# row before writing the values - otherwise Excel would silently
# convert "2024-10-03" to a date serial and "5"/"14" to numbers.
$newRows = @(
    @("2024-10-03", "5",  "uhf", "Heart Rate",     "VeryHigh", "Static"),
    @("2024-10-04", "14", "fr",  "Heart Rate",     "Low",      "Static"),
    @("2024-10-04", "14", "fr",  "Oxygen",         "High",     "Static"),
    @("2024-10-04", "14", "fr",  "Blood Pressure", "VeryHigh", "Decreasing")
)

$startRow = 9
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]

    $rowRange = $ws.Range("A" + $r + ":F" + $r)
    $rowRange.NumberFormat = "@"

    $ws.Range("A" + $r).Value = $rowValues[0]
    $ws.Range("B" + $r).Value = $rowValues[1]
    $ws.Range("C" + $r).Value = $rowValues[2]
    $ws.Range("D" + $r).Value = $rowValues[3]
    $ws.Range("E" + $r).Value = $rowValues[4]
    $ws.Range("F" + $r).Value = $rowValues[5]
}
